# Apply updated cryptocurrency price/volume data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'66.990.02"
$ws.Cells.Item(2, 5).Value = '  +0.62%  '

# Row 3
$ws.Cells.Item(3, 4).Value = "'3.500.38"
$ws.Cells.Item(3, 5).Value = '  -0.12%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.02%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'594.06"
$ws.Cells.Item(5, 5).Value = '  +0.51%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'173.10"
$ws.Cells.Item(6, 5).Value = '  +1.95%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.01%  '

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.604"
$ws.Cells.Item(8, 5).Value = '  +2.41%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +3.76%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -0.94%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.433"
$ws.Cells.Item(11, 5).Value = '  -1.31%  '

# Row 12
$ws.Cells.Item(12, 4).Value = "'4.104.82"
$ws.Cells.Item(12, 5).Value = '  -0.16%  '

# Row 13
$ws.Cells.Item(13, 4).Value = "'0.134"
$ws.Cells.Item(13, 5).Value = '  -0.24%  '

# Row 14
$ws.Cells.Item(14, 4).Value = "'28.91"
$ws.Cells.Item(14, 5).Value = '  +1.76%  '

# Row 15
$ws.Cells.Item(15, 4).Value = "'66.955.23"
$ws.Cells.Item(15, 5).Value = '  +0.52%  '

# Row 16
$ws.Cells.Item(16, 4).Value = "'0.0000178"
$ws.Cells.Item(16, 5).Value = '  -0.20%  '

# Row 17
$ws.Cells.Item(17, 4).Value = "'3.533.87"
$ws.Cells.Item(17, 5).Value = '  +1.18%  '

# Row 18
$ws.Cells.Item(18, 4).Value = "'6.31"
$ws.Cells.Item(18, 5).Value = '  -0.61%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.41%  '

# Row 20
$ws.Cells.Item(20, 4).Value = "'393.64"
$ws.Cells.Item(20, 5).Value = '  +0.76%  '

# Row 21
$ws.Cells.Item(21, 4).Value = "'7.99"
$ws.Cells.Item(21, 5).Value = '  -0.05%  '

# Row 22
$ws.Cells.Item(22, 4).Value = "'73.04"

# Row 23
$ws.Cells.Item(23, 5).Value = '  -0.07%  '

# Row 24
$ws.Cells.Item(24, 4).Value = "'0.536"
$ws.Cells.Item(24, 5).Value = '  -0.23%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -2.97%  '

# Row 26
$ws.Cells.Item(26, 4).Value = "'0.0000120"
$ws.Cells.Item(26, 5).Value = '  -2.09%  '

# Row 27
$ws.Cells.Item(27, 4).Value = "'10.19"
$ws.Cells.Item(27, 5).Value = '  -2.16%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +0.42%  '

# Row 29
$ws.Cells.Item(29, 4).Value = "'0.997"
$ws.Cells.Item(29, 5).Value = '  -0.37%  '

# Row 30
$ws.Cells.Item(30, 4).Value = "'6.26"
$ws.Cells.Item(30, 5).Value = '  -1.55%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -3.11%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -0.56%  '

# Row 33
$ws.Cells.Item(33, 4).Value = "'23.73"
$ws.Cells.Item(33, 5).Value = '  +0.40%  '

# Row 34
$ws.Cells.Item(34, 4).Value = "'7.34"
$ws.Cells.Item(34, 5).Value = '  -0.79%  '

# Row 35
$ws.Cells.Item(35, 4).Value = "'1.66"
$ws.Cells.Item(35, 5).Value = '  +2.31%  '

# Row 36
$ws.Cells.Item(36, 4).Value = "'163.44"
$ws.Cells.Item(36, 5).Value = '  +0.60%  '

# Row 37
$ws.Cells.Item(37, 4).Value = "'0.881"
$ws.Cells.Item(37, 5).Value = '  -0.24%  '

# Row 38
$ws.Cells.Item(38, 4).Value = "'1.89"
$ws.Cells.Item(38, 5).Value = '  -0.61%  '

# Row 39
$ws.Cells.Item(39, 4).Value = "'6.98"
$ws.Cells.Item(39, 5).Value = '  +2.49%  '

# Row 40
$ws.Cells.Item(40, 4).Value = "'4.67"
$ws.Cells.Item(40, 5).Value = '  -0.81%  '

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.0746"
$ws.Cells.Item(41, 5).Value = '  -0.19%  '

# Row 42
$ws.Cells.Item(42, 4).Value = "'27.32"
$ws.Cells.Item(42, 5).Value = '  -1.85%  '

# Row 43
$ws.Cells.Item(43, 4).Value = "'26.30"
$ws.Cells.Item(43, 5).Value = '  -0.88%  '

# Row 44
$ws.Cells.Item(44, 4).Value = "'2.807.14"
$ws.Cells.Item(44, 5).Value = '  -0.10%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'OKB'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(45, 4).Value = "'42.71"
$ws.Cells.Item(45, 5).Value = '  -1.07%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'dogwifhat'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(46, 4).Value = "'2.55"
$ws.Cells.Item(46, 5).Value = '  +0.38%  '

# Row 47
$ws.Cells.Item(47, 4).Value = "'0.0302"
$ws.Cells.Item(47, 5).Value = '  -3.02%  '

# Row 48
$ws.Cells.Item(48, 4).Value = "'335.59"
$ws.Cells.Item(48, 5).Value = '  -5.25%  '

# Row 49
$ws.Cells.Item(49, 4).Value = "'34.64"
$ws.Cells.Item(49, 5).Value = '  +2.50%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -0.97%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Stellar'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(51, 4).Value = "'0.106"
$ws.Cells.Item(51, 5).Value = '  +0.50%  '
